# Update the WIP sheet: rename the "COP 2025" campaign group (rows 3-8)
# to "COP 2026", and leave the grid selection on that range, matching
# the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WIP")
$ws.Activate()

$ws.Range("A3:A8").Value = "COP 2026"

$ws.Range("A3:A8").Select()
